# "task status sheet updated"
# The Task Status sheet (named "Sheet1") has a table where column D holds the
# person a module's task is assigned to. Update the "Assigned To" values for
# the first two modules so that the names are indented with leading spaces,
# matching the formatting already used further down the table (e.g. row 27).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1) Home Module -> Assigned To (merged D3:D5) : "Ayushi" -> "          Ayushi"
$ws.Range("D3").Value = "          Ayushi"

# 2) Product module -> Assigned To (merged D7:D8) : "Mayank" -> "        Mayank"
$ws.Range("D7").Value = "        Mayank"

# Update the on-screen selection to match where the edit left the cursor.
$ws.Activate()
$ws.Range("K11").Select()
